# Update cryptos list (GitHub Actions style refresh of D/E columns, plus a
# Polygon/Dogecoin row swap) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $value) {
    # Force text storage so Excel doesn't silently coerce values like
    # "291.80" or "49.89" into numbers (which would drop trailing zeros /
    # change the stored representation away from the scraped text).
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = "  $value  "
}

# Row 2 - Bitcoin
Set-Price 2 "22.456.53"
Set-Volume 2 "+0.33%"

# Row 3 - Ethereum
Set-Price 3 "1.574.17"
Set-Volume 3 "+0.19%"

# Row 4 - TetherUSD
Set-Volume 4 "-0.09%"

# Row 5 - USDC
Set-Volume 5 "-0.07%"

# Row 6 - BNB
Set-Price 6 "291.80"
Set-Volume 6 "+0.38%"

# Row 7 - XRP
Set-Price 7 "0.3730"

# Row 8 - OKB
Set-Price 8 "49.89"
Set-Volume 8 "-0.28%"

# Row 9 - Cardano
Set-Price 9 "0.3398"
Set-Volume 9 "-0.67%"

# Row 10 & 11 - Dogecoin / Polygon swap places
$ws.Cells.Item(10, 2).Value = "Polygon"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-Price 10 "1.145"
Set-Volume 10 "-0.45%"

$ws.Cells.Item(11, 2).Value = "Dogecoin"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-Price 11 "0.07567"
Set-Volume 11 "-0.80%"

# Row 12 - BinanceUSD
Set-Price 12 "1.002"
Set-Volume 12 "-0.12%"

# Row 13 - Solana
Set-Volume 13 "+0.77%"

# Row 14 - Polkadot
Set-Price 14 "6.022"
Set-Volume 14 "+0.01%"

# Row 15 - Chainlink
Set-Price 15 "6.964"
Set-Volume 15 "+0.54%"

# Row 16 - WrappedEther
Set-Price 16 "1.570.28"
Set-Volume 16 "-0.33%"

# Row 17 - ShibaInu
Set-Volume 17 "-0.41%"

# Row 18 - Litecoin
Set-Price 18 "90.94"
Set-Volume 18 "+1.03%"

# Row 19 - TRON
Set-Volume 19 "+0.26%"

# Row 20 - Dai
Set-Volume 20 "+0.04%"

# Row 21 - Uniswap
Set-Price 21 "6.305"
Set-Volume 21 "+1.87%"

# Row 22 - Avalanche
Set-Price 22 "16.32"
Set-Volume 22 "-2.77%"

# Row 23 - Cosmos
Set-Price 23 "12.17"
Set-Volume 23 "+1.59%"

# Row 24 - WrappedBTC
Set-Price 24 "22.459.80"
Set-Volume 24 "+0.37%"

# Row 25 - Toncoin
Set-Price 25 "2.341"
Set-Volume 25 "-2.39%"

# Row 26 - LidoDAOToken
Set-Price 26 "2.692"
Set-Volume 26 "+0.61%"

# Row 27 - EthereumClassic
Set-Price 27 "20.09"
Set-Volume 27 "-0.50%"

# Row 28 - Monero
Set-Price 28 "148.63"
Set-Volume 28 "+0.99%"

# Row 29 - HuobiToken
Set-Price 29 "5.004"
Set-Volume 29 "-0.59%"

# Row 30 - BitcoinCash
Set-Price 30 "125.57"
Set-Volume 30 "-0.44%"

# Row 31 - WrappedliquidstakedEther2.0
Set-Price 31 "1.747.88"
Set-Volume 31 "-0.03%"

# Row 32 - ImmutableX
Set-Price 32 "1.060"
Set-Volume 32 "+7.84%"

# Row 33 - Filecoin
Set-Price 33 "6.205"
Set-Volume 33 "+1.11%"

# Row 34 - WEMIXTOKEN
Set-Price 34 "1.986"
Set-Volume 34 "-1.03%"

# Row 35 - FraxShare
Set-Price 35 "9.834"
Set-Volume 35 "-0.38%"

# Row 36 - Stellar
Set-Price 36 "0.08384"
Set-Volume 36 "-1.69%"

# Row 37 - TrustWalletToken
Set-Price 37 "1.376"
Set-Volume 37 "+1.72%"

# Row 38 - VeChain
Set-Price 38 "0.02492"
Set-Volume 38 "-1.98%"

# Row 39 - Algorand
Set-Price 39 "0.2299"
Set-Volume 39 "-0.66%"

# Row 40 - Hedera
Set-Price 40 "0.06526"
Set-Volume 40 "-0.21%"

# Row 41 - InternetComputer(DFINITY)
Set-Price 41 "5.487"
Set-Volume 41 "+1.71%"

# Row 42 - Aptos
Set-Price 42 "11.30"
Set-Volume 42 "-0.86%"

# Row 43 - TheSandbox
Set-Price 43 "0.6235"
Set-Volume 43 "-2.03%"

# Row 44 - Frax
Set-Volume 44 "+0.00%"

# Row 45 - EnergySwap
Set-Price 45 "13.96"

# Row 46 - PancakeSwap
Set-Price 46 "3.815"
Set-Volume 46 "+0.87%"

# Row 47 - Decentraland
Set-Price 47 "0.5815"
Set-Volume 47 "-2.59%"

# Row 48 - Quant
Set-Price 48 "129.69"
Set-Volume 48 "+3.68%"

# Row 49 - NEARProtocol
Set-Price 49 "2.079"
Set-Volume 49 "-0.08%"

# Row 50 - EOS
Set-Price 50 "1.225"
Set-Volume 50 "-5.70%"

# Row 51 - Cronos
Set-Price 51 "0.07334"
Set-Volume 51 "+0.11%"
